$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 264
$ws.Range("I2").Value = 719
$ws.Range("J2").Value = 3108
$ws.Range("K2").Value = 18
$ws.Range("L2").Value = 821
$ws.Range("M2").Value = 47
$ws.Range("N2").Value = 524
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 41
$ws.Range("S2").Value = 315
$ws.Range("T2").Value = 489
$ws.Range("V2").Value = 4651
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 4670
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 69
$ws.Range("AA2").Value = 31
